$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = "'68.263.41"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.00%  '
$ws.Range('D3').Formula = "'2.450.41"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Formula = "'559.13"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.94%  '
$ws.Range('D6').Formula = "'163.17"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.32%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -2.10%  '
$ws.Range('D9').Formula = "'2.449.20"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.13%  '
$ws.Range('E10').Value = '  -5.93%  '
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('D12').Formula = "'0.338"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.60%  '
$ws.Range('D13').Formula = "'4.81"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.75%  '
$ws.Range('D14').Formula = "'2.906.78"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('D15').Formula = "'68.241.86"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.88%  '
$ws.Range('E16').Value = '  -3.92%  '
$ws.Range('D17').Formula = "'23.25"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.08%  '
$ws.Range('D18').Formula = "'2.464.01"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('D19').Formula = "'10.96"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('D20').Formula = "'7.16"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.26%  '
$ws.Range('D21').Formula = "'342.06"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.89%  '
$ws.Range('D22').Formula = "'3.78"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.48%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  -4.34%  '
$ws.Range('D25').Formula = "'67.73"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.08%  '
$ws.Range('D26').Formula = "'1.08"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.90%  '
$ws.Range('E27').Value = '  -6.44%  '
$ws.Range('D28').Formula = "'2.589.49"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('D29').Formula = "'8.15"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.23%  '
$ws.Range('D30').Formula = "'0.0₃0833"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.67%  '
$ws.Range('D31').Formula = "'7.23"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.67%  '
$ws.Range('D32').Formula = "'3.39"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +128.63%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Formula = "'432.19"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.37%  '
$ws.Range('E35').Value = '  -3.64%  '
$ws.Range('D36').Formula = "'1.67"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.71%  '
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('D38').Formula = "'18.99"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').Formula = "'0.109"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.94%  '
$ws.Range('D41').Formula = "'17.84"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.33%  '
$ws.Range('E42').Value = '  -3.64%  '
$ws.Range('D43').Formula = "'4.44"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.23%  '
$ws.Range('E44').Value = '  -5.34%  '
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('E46').Value = '  -6.54%  '
$ws.Range('D47').Formula = "'133.69"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.77%  '
$ws.Range('E48').Value = '  -3.74%  '
$ws.Range('E49').Value = '  -2.03%  '
$ws.Range('D50').Formula = "'0.483"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.15%  '
$ws.Range('D51').Formula = "'0.559"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.35%  '
